$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New ranking data is pasted in two column-wise passes (all of column A for
# rows 2-11, then all of column B for rows 2-11) in the same "raw" order the
# source data arrived in, before being sorted by Total Score (column E)
# descending. Writing the values in this exact order reproduces the shared
# string table layout of the target workbook.
# ---------------------------------------------------------------------------

# Column A (MT4 Account) - raw entry order
$ws.Range("A2").Value = "2xxx105"
$ws.Range("A4").Value = "2xxx485"
$ws.Range("A3").Value = "2xxx858"
$ws.Range("A5").Value = "2xxx556"
$ws.Range("A8").Value = "2xxx475"
$ws.Range("A6").Value = "2xxx647"
$ws.Range("A7").Value = "2xxx441"
$ws.Range("A10").Value = "2xxx213"
$ws.Range("A9").Value = "2xxx845"
$ws.Range("A11").Value = "2xxx447"

# Column B (Email) - same raw entry order
$ws.Range("B2").Value = "877******qq.com"
$ws.Range("B4").Value = "290*******qq.com"
$ws.Range("B3").Value = "she**************163.com"
$ws.Range("B5").Value = "Zfa*****************gmail.com"
$ws.Range("B8").Value = "774******qq.com"
$ws.Range("B6").Value = "skt******outlook.com"
$ws.Range("B7").Value = "233*******qq.com"
$ws.Range("B10").Value = "wan**************126.com"
$ws.Range("B9").Value = "fei*****163.com"
$ws.Range("B11").Value = "you*****qq.com"

# Make sure the new score cells carry the same number format as the existing
# score column (2 decimal places), matching style index already used by C2:E3
$numFmt = $ws.Range("C2").NumberFormat
$ws.Range("C2:E11").NumberFormat = $numFmt

# Column C (Net Equity Score) and D (Trading Volume Score)
$ws.Range("C2").Value = 35
$ws.Range("D2").Value = 5

$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 6

$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 4.3

$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 1

$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 1

$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 6.1

$ws.Range("C8").Value = 15
$ws.Range("D8").Value = 3

$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 1.4

$ws.Range("C10").Value = 13.33
$ws.Range("D10").Value = 5.2

$ws.Range("C11").Value = 11
$ws.Range("D11").Value = 0.5

# Column E (Total Score) - weighted formula, entered per-cell so each row
# keeps its own (non-shared) formula
$ws.Range("E2").Formula = "=C2*0.9+D2*0.1"
$ws.Range("E3").Formula = "=C3*0.9+D3*0.1"
$ws.Range("E4").Formula = "=C4*0.9+D4*0.1"
$ws.Range("E5").Formula = "=C5*0.9+D5*0.1"
$ws.Range("E6").Formula = "=C6*0.9+D6*0.1"
$ws.Range("E7").Formula = "=C7*0.9+D7*0.1"
$ws.Range("E8").Formula = "=C8*0.9+D8*0.1"
$ws.Range("E9").Formula = "=C9*0.9+D9*0.1"
$ws.Range("E10").Formula = "=C10*0.9+D10*0.1"
$ws.Range("E11").Formula = "=C11*0.9+D11*0.1"

# Re-apply the existing descending sort on Total Score across the now-larger
# data range (rows are already in the right order, so this just refreshes
# the workbook's sort state to cover A2:E11 without reshuffling anything).
$ws.Range("A1:E11").Sort($ws.Range("E1"), 2, $null, $null, 1, $null, 1, 1)
